$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their text formatting so numeric-looking
# strings (e.g. '1.00', '0.999') are not auto-converted to numbers by Excel.
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '42.690.53'
$ws.Range('E2').Value = '  -0.65%  '
$ws.Range('D3').Value = '2.509.01'
$ws.Range('E3').Value = '  -1.68%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '311.95'
$ws.Range('E5').Value = '  +2.49%  '
$ws.Range('D6').Value = '95.13'
$ws.Range('E6').Value = '  -3.30%  '
$ws.Range('D7').Value = '0.577'
$ws.Range('E7').Value = '  +0.33%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = '0.536'
$ws.Range('E9').Value = '  -1.73%  '
$ws.Range('D10').Value = '35.94'
$ws.Range('E10').Value = '  -1.16%  '
$ws.Range('D11').Value = '0.0808'
$ws.Range('E11').Value = '  -2.45%  '
$ws.Range('D12').Value = '7.67'
$ws.Range('E12').Value = '  -0.89%  '
$ws.Range('E13').Value = '  -2.08%  '
$ws.Range('D14').Value = '2.893.62'
$ws.Range('E14').Value = '  -1.69%  '
$ws.Range('D15').Value = '15.56'
$ws.Range('E15').Value = '  +4.37%  '
$ws.Range('D16').Value = '2.529.75'
$ws.Range('E16').Value = '  +1.49%  '
$ws.Range('D17').Value = '0.855'
$ws.Range('E17').Value = '  -2.06%  '
$ws.Range('D18').Value = '42.680.85'
$ws.Range('E18').Value = '  -0.77%  '
$ws.Range('D19').Value = '13.16'
$ws.Range('E19').Value = '  -2.65%  '
$ws.Range('D20').Value = '0.0₃0964'
$ws.Range('E20').Value = '  -2.98%  '
$ws.Range('D21').Value = '6.53'
$ws.Range('E21').Value = '  -1.24%  '
$ws.Range('D22').Value = '71.19'
$ws.Range('E22').Value = '  -1.12%  '
$ws.Range('D23').Value = '251.28'
$ws.Range('E23').Value = '  -1.27%  '
$ws.Range('D24').Value = '2.94'
$ws.Range('E24').Value = '  -0.62%  '
$ws.Range('E25').Value = '  -2.54%  '
$ws.Range('D26').Value = '26.71'
$ws.Range('E26').Value = '  -4.80%  '
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').Value = '2.35'
$ws.Range('E28').Value = '  +11.70%  '
$ws.Range('D29').Value = '39.06'
$ws.Range('E29').Value = '  +3.30%  '
$ws.Range('D30').Value = '10.02'
$ws.Range('E30').Value = '  -1.18%  '
$ws.Range('D31').Value = '5.89'
$ws.Range('E31').Value = '  -2.83%  '
$ws.Range('D32').Value = '156.69'
$ws.Range('E32').Value = '  -1.29%  '
$ws.Range('E33').Value = '  +3.61%  '
$ws.Range('D34').Value = '3.31'
$ws.Range('E34').Value = '  +0.17%  '
$ws.Range('D35').Value = '2.07'
$ws.Range('E35').Value = '  -4.03%  '
$ws.Range('D36').Value = '0.0784'
$ws.Range('E36').Value = '  -2.69%  '
$ws.Range('E37').Value = '  -5.45%  '
$ws.Range('E38').Value = '  -2.27%  '
$ws.Range('B39').Value = 'EnergySwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D39').Value = '24.20'
$ws.Range('E39').Value = '  -6.71%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').Value = '0.119'
$ws.Range('E40').Value = '  -0.31%  '
$ws.Range('D41').Value = '2.09'
$ws.Range('E41').Value = '  -0.13%  '
$ws.Range('D42').Value = '3.84'
$ws.Range('E42').Value = '  -1.67%  '
$ws.Range('D43').Value = '3.36'
$ws.Range('E43').Value = '  -2.07%  '
$ws.Range('D44').Value = '2.070.14'
$ws.Range('E44').Value = '  -0.77%  '
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('D46').Value = '0.0300'
$ws.Range('E46').Value = '  -2.01%  '
$ws.Range('D47').Value = '86.29'
$ws.Range('E47').Value = '  -0.29%  '
$ws.Range('E48').Value = '  -2.03%  '
$ws.Range('D49').Value = '2.752.43'
$ws.Range('E49').Value = '  -1.71%  '
$ws.Range('D50').Value = '73.23'
$ws.Range('E51').Value = '  -0.33%  '
